$d = $word.ActiveDocument

$d.Content.Find.Execute("0+83=83", $true, $false, $false, $false, $false, $true, 1, $false, "87-41=46", 2) | Out-Null
$d.Content.Find.Execute("26+57=83", $true, $false, $false, $false, $false, $true, 1, $false, "86-75=11", 2) | Out-Null
$d.Content.Find.Execute("62+36=98", $true, $false, $false, $false, $false, $true, 1, $false, "59-43=16", 2) | Out-Null
$d.Content.Find.Execute("26+37=63", $true, $false, $false, $false, $false, $true, 1, $false, "63-57=6", 2) | Out-Null
$d.Content.Find.Execute("46-15=31", $true, $false, $false, $false, $false, $true, 1, $false, "72-14=58", 2) | Out-Null
$d.Content.Find.Execute("62-58=4", $true, $false, $false, $false, $false, $true, 1, $false, "57-29=28", 2) | Out-Null
$d.Content.Find.Execute("11+13=24", $true, $false, $false, $false, $false, $true, 1, $false, "32-4=28", 2) | Out-Null
$d.Content.Find.Execute("67-8=59", $true, $false, $false, $false, $false, $true, 1, $false, "33+42=75", 2) | Out-Null
$d.Content.Find.Execute("16+83=99", $true, $false, $false, $false, $false, $true, 1, $false, "2+39=41", 2) | Out-Null
$d.Content.Find.Execute("5+9=14", $true, $false, $false, $false, $false, $true, 1, $false, "32+38=70", 2) | Out-Null
$d.Content.Find.Execute("4+17=21", $true, $false, $false, $false, $false, $true, 1, $false, "59-42=17", 2) | Out-Null
$d.Content.Find.Execute("97-34=63", $true, $false, $false, $false, $false, $true, 1, $false, "17+45=62", 2) | Out-Null
$d.Content.Find.Execute("83-35=48", $true, $false, $false, $false, $false, $true, 1, $false, "58-15=43", 2) | Out-Null
$d.Content.Find.Execute("10+51=61", $true, $false, $false, $false, $false, $true, 1, $false, "68-48=20", 2) | Out-Null
$d.Content.Find.Execute("15-12=3", $true, $false, $false, $false, $false, $true, 1, $false, "29+47=76", 2) | Out-Null
$d.Content.Find.Execute("88-34=54", $true, $false, $false, $false, $false, $true, 1, $false, "29+29=58", 2) | Out-Null
$d.Content.Find.Execute("75-56=19", $true, $false, $false, $false, $false, $true, 1, $false, "74-29=45", 2) | Out-Null
$d.Content.Find.Execute("49+10=59", $true, $false, $false, $false, $false, $true, 1, $false, "25+47=72", 2) | Out-Null
$d.Content.Find.Execute("47+12=59", $true, $false, $false, $false, $false, $true, 1, $false, "31+46=77", 2) | Out-Null
$d.Content.Find.Execute("66-18=48", $true, $false, $false, $false, $false, $true, 1, $false, "5+12=17", 2) | Out-Null
$d.Content.Find.Execute("62-60=2", $true, $false, $false, $false, $false, $true, 1, $false, "1+26=27", 2) | Out-Null
$d.Content.Find.Execute("87-30=57", $true, $false, $false, $false, $false, $true, 1, $false, "73+20=93", 2) | Out-Null
$d.Content.Find.Execute("23+18=41", $true, $false, $false, $false, $false, $true, 1, $false, "55-15=40", 2) | Out-Null
$d.Content.Find.Execute("26+26=52", $true, $false, $false, $false, $false, $true, 1, $false, "67-44=23", 2) | Out-Null
$d.Content.Find.Execute("72+20=92", $true, $false, $false, $false, $false, $true, 1, $false, "27-26=1", 2) | Out-Null
$d.Content.Find.Execute("22+69=91", $true, $false, $false, $false, $false, $true, 1, $false, "75-62=13", 2) | Out-Null
$d.Content.Find.Execute("42+13=55", $true, $false, $false, $false, $false, $true, 1, $false, "84-24=60", 2) | Out-Null
$d.Content.Find.Execute("22+26=48", $true, $false, $false, $false, $false, $true, 1, $false, "41+6=47", 2) | Out-Null
$d.Content.Find.Execute("68-53=15", $true, $false, $false, $false, $false, $true, 1, $false, "90+9=99", 2) | Out-Null
$d.Content.Find.Execute("20+37=57", $true, $false, $false, $false, $false, $true, 1, $false, "49-11=38", 2) | Out-Null
$d.Content.Find.Execute("87-78=9", $true, $false, $false, $false, $false, $true, 1, $false, "98-82=16", 2) | Out-Null
$d.Content.Find.Execute("66-17=49", $true, $false, $false, $false, $false, $true, 1, $false, "20+30=50", 2) | Out-Null
$d.Content.Find.Execute("72-3=69", $true, $false, $false, $false, $false, $true, 1, $false, "78-13=65", 2) | Out-Null
$d.Content.Find.Execute("60-21=39", $true, $false, $false, $false, $false, $true, 1, $false, "25+36=61", 2) | Out-Null
$d.Content.Find.Execute("63-62=1", $true, $false, $false, $false, $false, $true, 1, $false, "74-66=8", 2) | Out-Null
$d.Content.Find.Execute("41+53=94", $true, $false, $false, $false, $false, $true, 1, $false, "62+20=82", 2) | Out-Null
$d.Content.Find.Execute("23+25=48", $true, $false, $false, $false, $false, $true, 1, $false, "94-94=0", 2) | Out-Null
$d.Content.Find.Execute("76-63=13", $true, $false, $false, $false, $false, $true, 1, $false, "19+29=48", 2) | Out-Null
$d.Content.Find.Execute("50+22=72", $true, $false, $false, $false, $false, $true, 1, $false, "60-35=25", 2) | Out-Null
$d.Content.Find.Execute("41-4=37", $true, $false, $false, $false, $false, $true, 1, $false, "24-9=15", 2) | Out-Null
$d.Content.Find.Execute("8-4=4", $true, $false, $false, $false, $false, $true, 1, $false, "57+24=81", 2) | Out-Null
$d.Content.Find.Execute("81-30=51", $true, $false, $false, $false, $false, $true, 1, $false, "95-8=87", 2) | Out-Null
$d.Content.Find.Execute("61-1=60", $true, $false, $false, $false, $false, $true, 1, $false, "94-43=51", 2) | Out-Null
$d.Content.Find.Execute("29-15=14", $true, $false, $false, $false, $false, $true, 1, $false, "89-45=44", 2) | Out-Null
$d.Content.Find.Execute("2+31=33", $true, $false, $false, $false, $false, $true, 1, $false, "92-43=49", 2) | Out-Null
$d.Content.Find.Execute("85+14=99", $true, $false, $false, $false, $false, $true, 1, $false, "74-50=24", 2) | Out-Null
$d.Content.Find.Execute("6+81=87", $true, $false, $false, $false, $false, $true, 1, $false, "23+0=23", 2) | Out-Null
$d.Content.Find.Execute("91-1=90", $true, $false, $false, $false, $false, $true, 1, $false, "88-88=0", 2) | Out-Null
$d.Content.Find.Execute("82-67=15", $true, $false, $false, $false, $false, $true, 1, $false, "79-12=67", 2) | Out-Null
$d.Content.Find.Execute("98-30=68", $true, $false, $false, $false, $false, $true, 1, $false, "72+25=97", 2) | Out-Null
$d.Content.Find.Execute("54-45=9", $true, $false, $false, $false, $false, $true, 1, $false, "78-73=5", 2) | Out-Null
$d.Content.Find.Execute("1+19=20", $true, $false, $false, $false, $false, $true, 1, $false, "35-5=30", 2) | Out-Null
$d.Content.Find.Execute("34+53=87", $true, $false, $false, $false, $false, $true, 1, $false, "97-55=42", 2) | Out-Null
$d.Content.Find.Execute("50+16=66", $true, $false, $false, $false, $false, $true, 1, $false, "96-33=63", 2) | Out-Null
$d.Content.Find.Execute("40+50=90", $true, $false, $false, $false, $false, $true, 1, $false, "1+51=52", 2) | Out-Null
$d.Content.Find.Execute("42-14=28", $true, $false, $false, $false, $false, $true, 1, $false, "15+70=85", 2) | Out-Null
$d.Content.Find.Execute("79-25=54", $true, $false, $false, $false, $false, $true, 1, $false, "69+16=85", 2) | Out-Null
$d.Content.Find.Execute("42+17=59", $true, $false, $false, $false, $false, $true, 1, $false, "4+8=12", 2) | Out-Null
$d.Content.Find.Execute("94-69=25", $true, $false, $false, $false, $false, $true, 1, $false, "54-23=31", 2) | Out-Null
$d.Content.Find.Execute("83-40=43", $true, $false, $false, $false, $false, $true, 1, $false, "92-29=63", 2) | Out-Null
$d.Content.Find.Execute("84-70=14", $true, $false, $false, $false, $false, $true, 1, $false, "68+11=79", 2) | Out-Null
$d.Content.Find.Execute("85-4=81", $true, $false, $false, $false, $false, $true, 1, $false, "94-0=94", 2) | Out-Null
$d.Content.Find.Execute("85-65=20", $true, $false, $false, $false, $false, $true, 1, $false, "89+3=92", 2) | Out-Null
$d.Content.Find.Execute("96-59=37", $true, $false, $false, $false, $false, $true, 1, $false, "34-2=32", 2) | Out-Null
$d.Content.Find.Execute("16+39=55", $true, $false, $false, $false, $false, $true, 1, $false, "34-15=19", 2) | Out-Null
$d.Content.Find.Execute("46+23=69", $true, $false, $false, $false, $false, $true, 1, $false, "93-92=1", 2) | Out-Null
$d.Content.Find.Execute("39+4=43", $true, $false, $false, $false, $false, $true, 1, $false, "61+15=76", 2) | Out-Null
$d.Content.Find.Execute("63+11=74", $true, $false, $false, $false, $false, $true, 1, $false, "54+28=82", 2) | Out-Null
$d.Content.Find.Execute("90-27=63", $true, $false, $false, $false, $false, $true, 1, $false, "35+20=55", 2) | Out-Null
$d.Content.Find.Execute("43+4=47", $true, $false, $false, $false, $false, $true, 1, $false, "33+42=75", 2) | Out-Null
$d.Content.Find.Execute("85-36=49", $true, $false, $false, $false, $false, $true, 1, $false, "17+77=94", 2) | Out-Null
$d.Content.Find.Execute("54+1=55", $true, $false, $false, $false, $false, $true, 1, $false, "19+23=42", 2) | Out-Null
$d.Content.Find.Execute("57+37=94", $true, $false, $false, $false, $false, $true, 1, $false, "62-53=9", 2) | Out-Null
$d.Content.Find.Execute("58-2=56", $true, $false, $false, $false, $false, $true, 1, $false, "26+0=26", 2) | Out-Null
$d.Content.Find.Execute("34-20=14", $true, $false, $false, $false, $false, $true, 1, $false, "88+9=97", 2) | Out-Null
$d.Content.Find.Execute("38-29=9", $true, $false, $false, $false, $false, $true, 1, $false, "29+48=77", 2) | Out-Null
$d.Content.Find.Execute("10+54=64", $true, $false, $false, $false, $false, $true, 1, $false, "91-52=39", 2) | Out-Null
$d.Content.Find.Execute("8+28=36", $true, $false, $false, $false, $false, $true, 1, $false, "19+57=76", 2) | Out-Null
$d.Content.Find.Execute("91-84=7", $true, $false, $false, $false, $false, $true, 1, $false, "74-62=12", 2) | Out-Null
$d.Content.Find.Execute("63-25=38", $true, $false, $false, $false, $false, $true, 1, $false, "39-28=11", 2) | Out-Null
$d.Content.Find.Execute("68-57=11", $true, $false, $false, $false, $false, $true, 1, $false, "48-24=24", 2) | Out-Null
$d.Content.Find.Execute("10-8=2", $true, $false, $false, $false, $false, $true, 1, $false, "15+83=98", 2) | Out-Null
$d.Content.Find.Execute("71-26=45", $true, $false, $false, $false, $false, $true, 1, $false, "15+22=37", 2) | Out-Null
$d.Content.Find.Execute("54+7=61", $true, $false, $false, $false, $false, $true, 1, $false, "34-16=18", 2) | Out-Null
$d.Content.Find.Execute("94-36=58", $true, $false, $false, $false, $false, $true, 1, $false, "77+4=81", 2) | Out-Null
$d.Content.Find.Execute("99-93=6", $true, $false, $false, $false, $false, $true, 1, $false, "39-38=1", 2) | Out-Null
$d.Content.Find.Execute("48-18=30", $true, $false, $false, $false, $false, $true, 1, $false, "23+49=72", 2) | Out-Null
$d.Content.Find.Execute("81-67=14", $true, $false, $false, $false, $false, $true, 1, $false, "13+42=55", 2) | Out-Null
$d.Content.Find.Execute("95-71=24", $true, $false, $false, $false, $false, $true, 1, $false, "86-56=30", 2) | Out-Null
$d.Content.Find.Execute("70+1=71", $true, $false, $false, $false, $false, $true, 1, $false, "14-3=11", 2) | Out-Null
$d.Content.Find.Execute("43-13=30", $true, $false, $false, $false, $false, $true, 1, $false, "96-30=66", 2) | Out-Null
$d.Content.Find.Execute("40+3=43", $true, $false, $false, $false, $false, $true, 1, $false, "42-32=10", 2) | Out-Null
$d.Content.Find.Execute("89-84=5", $true, $false, $false, $false, $false, $true, 1, $false, "18+26=44", 2) | Out-Null
$d.Content.Find.Execute("99-5=94", $true, $false, $false, $false, $false, $true, 1, $false, "32-14=18", 2) | Out-Null
$d.Content.Find.Execute("11+30=41", $true, $false, $false, $false, $false, $true, 1, $false, "68-13=55", 2) | Out-Null
$d.Content.Find.Execute("57-3=54", $true, $false, $false, $false, $false, $true, 1, $false, "94-72=22", 2) | Out-Null
$d.Content.Find.Execute("15-10=5", $true, $false, $false, $false, $false, $true, 1, $false, "77-1=76", 2) | Out-Null
$d.Content.Find.Execute("29+14=43", $true, $false, $false, $false, $false, $true, 1, $false, "34+0=34", 2) | Out-Null
$d.Content.Find.Execute("93-28=65", $true, $false, $false, $false, $false, $true, 1, $false, "87-87=0", 2) | Out-Null
$d.Content.Find.Execute("40+6=46", $true, $false, $false, $false, $false, $true, 1, $false, "42-26=16", 2) | Out-Null
